# LodeRunner v3 - extend FrameCompare data up to Lv20 (commit: "lode runner v3 - up to lv20")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3")
$ws.Activate()

# Newly recorded level start/end frame markers (Lv13 End .. Lv20 Start).
# Column A = label, B = v3 frame, C = v2 frame, D (formula already present) = C - B
$data = @(
    @(27, "Lv13 End",   16257, 45855),
    @(28, "Lv14 Start", 16601, 46921),
    @(29, "Lv14 End",   17482, 50630),
    @(30, "Lv15 Start", 17826, 51696),
    @(31, "Lv15 End",   18918, 55302),
    @(32, "Lv16 Start", 19261, 56368),
    @(33, "Lv16 End",   19995, 58398),
    @(34, "Lv17 Start", 20338, 59464),
    @(35, "Lv17 End",   21653, 63242),
    @(36, "Lv18 Start", 21996, 64308),
    @(37, "Lv18 End",   23337, 68785),
    @(38, "Lv19 Start", 23681, 69851),
    @(39, "Lv19 End",   24335, 71616),
    @(40, "Lv20 Start", 24677, 72682)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Update the saved view state: scrolled position and active selection.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select()
